$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column width adjustments (C widened, T narrowed; engine quantizes
# ColumnWidth to a 1/6-character pixel grid, so the closest achievable
# values are used).
$ws.Range("C1").ColumnWidth = 22.833333333333336
$ws.Range("T1").ColumnWidth = 21.833333333333336

# Updated ratio/error values (recalculated error figures)
$ws.Range("B2").Value = -1.576667569721746
$ws.Range("C2").Value = 0.001478752354046088
$ws.Range("F2").Value = 0.007144296897712548
$ws.Range("G2").Value = 0.007417903189274322
$ws.Range("J2").Value = 2.083253286508519
$ws.Range("K2").Value = 0.6590415179803142
$ws.Range("L2").Value = 0.007567191028275208
$ws.Range("M2").Value = 0.1481087536733177
$ws.Range("N2").Value = 0.00005488204341624448
$ws.Range("O2").Value = 0.1481087536733174
$ws.Range("P2").Value = 1.325488815396536
$ws.Range("Q2").Value = 0.3215822963079741
$ws.Range("T2").Value = 0.2425649162228197
$ws.Range("U2").Value = 0.3269535019506545
$ws.Range("H3").Value = 5.185689931824315
$ws.Range("I3").Value = 0.005759621605348021
$ws.Range("P3").Value = 0.003275356462361479
$ws.Range("Q3").Value = 2.119924914268915
$ws.Range("R3").Value = 0.09733120204652208
$ws.Range("S3").Value = 0.2060866016614308
$ws.Range("T3").Value = 0.0003203543385384721
$ws.Range("U3").Value = 2.020137790100577
$ws.Range("B4").Value = -1.993968138692681
$ws.Range("C4").Value = 0.001741884743836049
$ws.Range("D4").Value = 0.0100170998366197
$ws.Range("E4").Value = 0.8474172615102259
$ws.Range("H4").Value = 2.732303424305919
$ws.Range("I4").Value = 0.01120578579880862
$ws.Range("J4").Value = 2.073822321355093
$ws.Range("K4").Value = 0.8201294338742942
$ws.Range("L4").Value = 0.007564028248501297
$ws.Range("M4").Value = 0.1745364945928616
$ws.Range("N4").Value = 0.00005485910494195209
$ws.Range("O4").Value = 0.1745364945928622
$ws.Range("P4").Value = 1.323289409733625
$ws.Range("Q4").Value = 0.2736660921300703
$ws.Range("R4").Value = 0.1830152333401183
$ws.Range("S4").Value = 0.3958132555503516
$ws.Range("T4").Value = 0.2414261194212968
$ws.Range("U4").Value = 0.3636864051450315
$ws.Range("P5").Value = 0.002818224514931116
$ws.Range("Q5").Value = 1.876895080002627
$ws.Range("T5").Value = 0.000512323621250908
$ws.Range("U5").Value = 1.970959586851783
$ws.Range("P6").Value = 1.32219823094449
$ws.Range("Q6").Value = 0.2209739728222805
$ws.Range("T6").Value = 0.2396666663667818
$ws.Range("U6").Value = 0.3656704898636055
$ws.Range("D7").Value = 0.009976845361061415
$ws.Range("E7").Value = 0.3297874133506366
$ws.Range("P7").Value = 0.003639453782854304
$ws.Range("Q7").Value = 1.822070429031099
$ws.Range("T7").Value = 0.0006632840911626782
$ws.Range("U7").Value = 1.940759177824897
$ws.Range("P8").Value = 1.319622257679337
$ws.Range("Q8").Value = 0.2325182917524557
$ws.Range("T8").Value = 0.2382819122323587
$ws.Range("U8").Value = 0.3670471312898169
$ws.Range("B9").Value = 144.2477879671713
$ws.Range("C9").Value = 0.0003612035633638828
$ws.Range("D9").Value = 0.00998226392451965
$ws.Range("E9").Value = 0.3367852707192005
$ws.Range("J9").Value = 4.698652581306738
$ws.Range("K9").Value = 0.3404243708834113
$ws.Range("L9").Value = 0.008672415110885429
$ws.Range("M9").Value = 0.03156690073271559
$ws.Range("N9").Value = 0.0000628978257402066
$ws.Range("O9").Value = 0.03156690073271488
$ws.Range("P9").Value = 0.003179041188988946
$ws.Range("Q9").Value = 2.1377126504212
$ws.Range("T9").Value = 0.0003620849407583208
$ws.Range("U9").Value = 2.044010697369161
$ws.Range("B10").Value = -1.507900421679476
$ws.Range("C10").Value = 0.001813883976860164
$ws.Range("L10").Value = 0.007567712224173581
$ws.Range("M10").Value = 0.1816623263845749
$ws.Range("N10").Value = 0.00005488582345771776
$ws.Range("O10").Value = 0.181662326384575
$ws.Range("P10").Value = 1.317214722714466
$ws.Range("Q10").Value = 0.2582171694122649
$ws.Range("T10").Value = 0.2367852110607417
$ws.Range("U10").Value = 0.4260308980406651
$ws.Range("P11").Value = 0.003396968616918185
$ws.Range("Q11").Value = 2.26423125468954
$ws.Range("R11").Value = 0.2061483249746852
$ws.Range("S11").Value = 0.3523134073429333
$ws.Range("T11").Value = 0.0007049510064866325
$ws.Range("U11").Value = 2.297564997300947
$ws.Range("F12").Value = 0.007145697660176736
$ws.Range("G12").Value = 0.006052072924899333
$ws.Range("P12").Value = 1.317266026380881
$ws.Range("Q12").Value = 0.245071663026405
$ws.Range("T12").Value = 0.2384620925124516
$ws.Range("U12").Value = 0.4561006561394234
